# Natmi following Dr Hou advice
# Update LR-pair table: add FAPs as a possible target cluster (rows now
# cover the full ECs/FAPs/sCs x ECs/FAPs/sCs combinations) and refresh stats.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs -> ECs
$ws.Cells.Item(2,1).Value = "ECs"
$ws.Cells.Item(2,2).Value = "Spp1"
$ws.Cells.Item(2,3).Value = "Itga4"
$ws.Cells.Item(2,4).Value = "ECs"
$ws.Cells.Item(2,5).Value = 3
$ws.Cells.Item(2,6).Value = 1
$ws.Cells.Item(2,7).Value = 282.78797966666667207392
$ws.Cells.Item(2,8).Value = 848.36393899999995937833
$ws.Cells.Item(2,9).Value = 0.96745217414012674162
$ws.Cells.Item(2,10).Value = 0.9674521741401266306
$ws.Cells.Item(2,11).Value = 3
$ws.Cells.Item(2,12).Value = 1
$ws.Cells.Item(2,13).Value = 26.04517333333333084511
$ws.Cells.Item(2,14).Value = 78.13551999999999964075
$ws.Cells.Item(2,15).Value = 0.9210237118384171362
$ws.Cells.Item(2,16).Value = 0.92102371183841702518
$ws.Cells.Item(2,17).Value = 7365.2619470014760736376
$ws.Cells.Item(2,18).Value = 66287.35752301327011082321
$ws.Cells.Item(2,19).Value = 0.89104639245268624226
$ws.Cells.Item(2,20).Value = 0.89104639245268602021

# Row 3: ECs -> FAPs
$ws.Cells.Item(3,1).Value = "ECs"
$ws.Cells.Item(3,2).Value = "Spp1"
$ws.Cells.Item(3,3).Value = "Itga4"
$ws.Cells.Item(3,4).Value = "FAPs"
$ws.Cells.Item(3,5).Value = 3
$ws.Cells.Item(3,6).Value = 1
$ws.Cells.Item(3,7).Value = 282.78797966666667207392
$ws.Cells.Item(3,8).Value = 848.36393899999995937833
$ws.Cells.Item(3,9).Value = 0.96745217414012674162
$ws.Cells.Item(3,10).Value = 0.9674521741401266306
$ws.Cells.Item(3,11).Value = 2
$ws.Cells.Item(3,12).Value = 0.66666666666666662966
$ws.Cells.Item(3,13).Value = 0.33022233333333328442
$ws.Cells.Item(3,14).Value = 0.99066699999999996429
$ws.Cells.Item(3,15).Value = 0.01167750336256581992
$ws.Cells.Item(3,16).Value = 0.01167750336256581992
$ws.Cells.Item(3,17).Value = 93.38290648414589156801
$ws.Cells.Item(3,18).Value = 840.44615835731292463606
$ws.Cells.Item(3,19).Value = 0.01129742601664295021
$ws.Cells.Item(3,20).Value = 0.01129742601664295021

# Row 4: ECs -> sCs
$ws.Cells.Item(4,1).Value = "ECs"
$ws.Cells.Item(4,2).Value = "Spp1"
$ws.Cells.Item(4,3).Value = "Itga4"
$ws.Cells.Item(4,4).Value = "sCs"
$ws.Cells.Item(4,5).Value = 3
$ws.Cells.Item(4,6).Value = 1
$ws.Cells.Item(4,7).Value = 282.78797966666667207392
$ws.Cells.Item(4,8).Value = 848.36393899999995937833
$ws.Cells.Item(4,9).Value = 0.96745217414012674162
$ws.Cells.Item(4,10).Value = 0.9674521741401266306
$ws.Cells.Item(4,11).Value = 3
$ws.Cells.Item(4,12).Value = 1
$ws.Cells.Item(4,13).Value = 1.90310899999999993959
$ws.Cells.Item(4,14).Value = 5.70932700000000004081
$ws.Cells.Item(4,15).Value = 0.06729878479901707511
$ws.Cells.Item(4,16).Value = 0.06729878479901707511
$ws.Cells.Item(4,17).Value = 538.17634919545037064381
$ws.Cells.Item(4,18).Value = 4843.58714275905276736012
$ws.Cells.Item(4,19).Value = 0.06510835567079757691
$ws.Cells.Item(4,20).Value = 0.06510835567079757691

# Row 5: FAPs -> ECs
$ws.Cells.Item(5,1).Value = "FAPs"
$ws.Cells.Item(5,2).Value = "Spp1"
$ws.Cells.Item(5,3).Value = "Itga4"
$ws.Cells.Item(5,4).Value = "ECs"
$ws.Cells.Item(5,5).Value = 3
$ws.Cells.Item(5,6).Value = 1
$ws.Cells.Item(5,7).Value = 7.71469666666666764598
$ws.Cells.Item(5,8).Value = 23.14408999999999849706
$ws.Cells.Item(5,9).Value = 0.02639291836872236993
$ws.Cells.Item(5,10).Value = 0.02639291836872236993
$ws.Cells.Item(5,11).Value = 3
$ws.Cells.Item(5,12).Value = 1
$ws.Cells.Item(5,13).Value = 26.04517333333333084511
$ws.Cells.Item(5,14).Value = 78.13551999999999964075
$ws.Cells.Item(5,15).Value = 0.9210237118384171362
$ws.Cells.Item(5,16).Value = 0.92102371183841702518
$ws.Cells.Item(5,17).Value = 200.93061189742229544208
$ws.Cells.Item(5,18).Value = 1808.37550707679997685773
$ws.Cells.Item(5,19).Value = 0.02430850364220902074
$ws.Cells.Item(5,20).Value = 0.02430850364220902074

# Row 6: FAPs -> FAPs
$ws.Cells.Item(6,1).Value = "FAPs"
$ws.Cells.Item(6,2).Value = "Spp1"
$ws.Cells.Item(6,3).Value = "Itga4"
$ws.Cells.Item(6,4).Value = "FAPs"
$ws.Cells.Item(6,5).Value = 3
$ws.Cells.Item(6,6).Value = 1
$ws.Cells.Item(6,7).Value = 7.71469666666666764598
$ws.Cells.Item(6,8).Value = 23.14408999999999849706
$ws.Cells.Item(6,9).Value = 0.02639291836872236993
$ws.Cells.Item(6,10).Value = 0.02639291836872236993
$ws.Cells.Item(6,11).Value = 2
$ws.Cells.Item(6,12).Value = 0.66666666666666662966
$ws.Cells.Item(6,13).Value = 0.33022233333333328442
$ws.Cells.Item(6,14).Value = 0.99066699999999996429
$ws.Cells.Item(6,15).Value = 0.01167750336256581992
$ws.Cells.Item(6,16).Value = 0.01167750336256581992
$ws.Cells.Item(6,17).Value = 2.54756513422555608273
$ws.Cells.Item(6,18).Value = 22.92808620802999897137
$ws.Cells.Item(6,19).Value = 0.0003082033929986808
$ws.Cells.Item(6,20).Value = 0.00030820339299868069

# Row 7: FAPs -> sCs
$ws.Cells.Item(7,1).Value = "FAPs"
$ws.Cells.Item(7,2).Value = "Spp1"
$ws.Cells.Item(7,3).Value = "Itga4"
$ws.Cells.Item(7,4).Value = "sCs"
$ws.Cells.Item(7,5).Value = 3
$ws.Cells.Item(7,6).Value = 1
$ws.Cells.Item(7,7).Value = 7.71469666666666764598
$ws.Cells.Item(7,8).Value = 23.14408999999999849706
$ws.Cells.Item(7,9).Value = 0.02639291836872236993
$ws.Cells.Item(7,10).Value = 0.02639291836872236993
$ws.Cells.Item(7,11).Value = 3
$ws.Cells.Item(7,12).Value = 1
$ws.Cells.Item(7,13).Value = 1.90310899999999993959
$ws.Cells.Item(7,14).Value = 5.70932700000000004081
$ws.Cells.Item(7,15).Value = 0.06729878479901707511
$ws.Cells.Item(7,16).Value = 0.06729878479901707511
$ws.Cells.Item(7,17).Value = 14.68190865860334071158
$ws.Cells.Item(7,18).Value = 132.13717792743000245537
$ws.Cells.Item(7,19).Value = 0.00177621133351467202
$ws.Cells.Item(7,20).Value = 0.00177621133351467202

# Row 8: sCs -> ECs
$ws.Cells.Item(8,1).Value = "sCs"
$ws.Cells.Item(8,2).Value = "Spp1"
$ws.Cells.Item(8,3).Value = "Itga4"
$ws.Cells.Item(8,4).Value = "ECs"
$ws.Cells.Item(8,5).Value = 3
$ws.Cells.Item(8,6).Value = 1
$ws.Cells.Item(8,7).Value = 1.79909033333333301385
$ws.Cells.Item(8,8).Value = 5.39727099999999992974
$ws.Cells.Item(8,9).Value = 0.00615490749115098299
$ws.Cells.Item(8,10).Value = 0.00615490749115098299
$ws.Cells.Item(8,11).Value = 3
$ws.Cells.Item(8,12).Value = 1
$ws.Cells.Item(8,13).Value = 26.04517333333333084511
$ws.Cells.Item(8,14).Value = 78.13551999999999964075
$ws.Cells.Item(8,15).Value = 0.9210237118384171362
$ws.Cells.Item(8,16).Value = 0.92102371183841702518
$ws.Cells.Item(8,17).Value = 46.85761957399110855249
$ws.Cells.Item(8,18).Value = 421.71857616592001249956
$ws.Cells.Item(8,19).Value = 0.0056688157435219582
$ws.Cells.Item(8,20).Value = 0.00566881574352195734

# Row 9: sCs -> FAPs
$ws.Cells.Item(9,1).Value = "sCs"
$ws.Cells.Item(9,2).Value = "Spp1"
$ws.Cells.Item(9,3).Value = "Itga4"
$ws.Cells.Item(9,4).Value = "FAPs"
$ws.Cells.Item(9,5).Value = 3
$ws.Cells.Item(9,6).Value = 1
$ws.Cells.Item(9,7).Value = 1.79909033333333301385
$ws.Cells.Item(9,8).Value = 5.39727099999999992974
$ws.Cells.Item(9,9).Value = 0.00615490749115098299
$ws.Cells.Item(9,10).Value = 0.00615490749115098299
$ws.Cells.Item(9,11).Value = 2
$ws.Cells.Item(9,12).Value = 0.66666666666666662966
$ws.Cells.Item(9,13).Value = 0.33022233333333328442
$ws.Cells.Item(9,14).Value = 0.99066699999999996429
$ws.Cells.Item(9,15).Value = 0.01167750336256581992
$ws.Cells.Item(9,16).Value = 0.01167750336256581992
$ws.Cells.Item(9,17).Value = 0.59409980775077775306
$ws.Cells.Item(9,18).Value = 5.34689826975699933342
$ws.Cells.Item(9,19).Value = 0.00007187395292419717
$ws.Cells.Item(9,20).Value = 0.00007187395292419717

# Row 10: sCs -> sCs
$ws.Cells.Item(10,1).Value = "sCs"
$ws.Cells.Item(10,2).Value = "Spp1"
$ws.Cells.Item(10,3).Value = "Itga4"
$ws.Cells.Item(10,4).Value = "sCs"
$ws.Cells.Item(10,5).Value = 3
$ws.Cells.Item(10,6).Value = 1
$ws.Cells.Item(10,7).Value = 1.79909033333333301385
$ws.Cells.Item(10,8).Value = 5.39727099999999992974
$ws.Cells.Item(10,9).Value = 0.00615490749115098299
$ws.Cells.Item(10,10).Value = 0.00615490749115098299
$ws.Cells.Item(10,11).Value = 3
$ws.Cells.Item(10,12).Value = 1
$ws.Cells.Item(10,13).Value = 1.90310899999999993959
$ws.Cells.Item(10,14).Value = 5.70932700000000004081
$ws.Cells.Item(10,15).Value = 0.06729878479901707511
$ws.Cells.Item(10,16).Value = 0.06729878479901707511
$ws.Cells.Item(10,17).Value = 3.42386500517966618773
$ws.Cells.Item(10,18).Value = 30.81478504661700057454
$ws.Cells.Item(10,19).Value = 0.00041421779470482808
$ws.Cells.Item(10,20).Value = 0.00041421779470482808
